# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the per-class Moogle_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR)
# per the scheduled-runner profit recompute.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2767.8235
$ws.Range("J17").Value = 2767.8235
$ws.Range("L17").Value = 8303.470499999999
$ws.Range("N17").Value = -8639.470499999999
$ws.Range("H19").Value = 52602.45
$ws.Range("I19").Value = 112896.555
$ws.Range("K19").Value = 112896.555
$ws.Range("M19").Value = -112721.555
$ws.Range("H96").Value = 763.4211
$ws.Range("I96").Value = 536.2222
$ws.Range("J96").Value = 967.9
$ws.Range("K96").Value = 1608.6666
$ws.Range("L96").Value = 2903.7
$ws.Range("M96").Value = -235.6666
$ws.Range("N96").Value = -5649.7
$ws.Range("H103").Value = 2232.0667
$ws.Range("I103").Value = 2298.6
$ws.Range("J103").Value = 2099
$ws.Range("K103").Value = 6895.799999999999
$ws.Range("L103").Value = 6297
$ws.Range("M103").Value = -6309.799999999999
$ws.Range("N103").Value = -7469
$ws.Range("H129").Value = 4180.1875
$ws.Range("I129").Value = 4319.6665
$ws.Range("J129").Value = 2088
$ws.Range("K129").Value = 12958.9995
$ws.Range("L129").Value = 6264
$ws.Range("M129").Value = -7958.999500000002
$ws.Range("N129").Value = -16264
$ws.Range("H135").Value = 2457.7273
$ws.Range("I135").Value = 1860.4
$ws.Range("K135").Value = 16743.6
$ws.Range("M135").Value = -14208.6
$ws.Range("H138").Value = 2990.8
$ws.Range("I138").Value = 2393.0278
$ws.Range("J138").Value = 3732.862
$ws.Range("K138").Value = 7179.0834
$ws.Range("L138").Value = 11198.586
$ws.Range("M138").Value = -2039.0834
$ws.Range("N138").Value = -21478.586
$ws.Range("H141").Value = 3542.75
$ws.Range("I141").Value = 1827.1923
$ws.Range("J141").Value = 6728.7856
$ws.Range("K141").Value = 5481.5769
$ws.Range("L141").Value = 20186.3568
$ws.Range("M141").Value = -301.5769
$ws.Range("N141").Value = -30546.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1372.95
$ws.Range("I2").Value = 1645.5
$ws.Range("K2").Value = 1645.5
$ws.Range("M2").Value = -1532.5
$ws.Range("H32").Value = 4401.032
$ws.Range("I32").Value = 2963.0715
$ws.Range("J32").Value = 17822
$ws.Range("K32").Value = 2963.0715
$ws.Range("L32").Value = 17822
$ws.Range("M32").Value = -2676.0715
$ws.Range("N32").Value = -18396
$ws.Range("H116").Value = 1372.95
$ws.Range("I116").Value = 1645.5
$ws.Range("K116").Value = 1645.5
$ws.Range("M116").Value = 648.5
$ws.Range("H122").Value = 2523.6155
$ws.Range("I122").Value = 1824.56
$ws.Range("K122").Value = 5473.68
$ws.Range("M122").Value = -3023.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1372.95
$ws.Range("I3").Value = 1645.5
$ws.Range("K3").Value = 1645.5
$ws.Range("M3").Value = -1531.5
$ws.Range("H20").Value = 3050.3948
$ws.Range("I20").Value = 3086.423
$ws.Range("K20").Value = 3086.423
$ws.Range("M20").Value = -2839.423
$ws.Range("H134").Value = 1202.1923
$ws.Range("I134").Value = 1206.32
$ws.Range("K134").Value = 3618.96
$ws.Range("M134").Value = -1083.96

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3924.132
$ws.Range("I58").Value = 2169.8647
$ws.Range("J58").Value = 7980.875
$ws.Range("K58").Value = 2169.8647
$ws.Range("L58").Value = 7980.875
$ws.Range("M58").Value = -1966.8647
$ws.Range("N58").Value = -8386.875
$ws.Range("H86").Value = 6798.8
$ws.Range("J86").Value = 7999.6665
$ws.Range("L86").Value = 7999.6665
$ws.Range("N86").Value = -10245.6665
$ws.Range("H89").Value = 6798.8
$ws.Range("J89").Value = 7999.6665
$ws.Range("L89").Value = 39998.3325
$ws.Range("N89").Value = -51230.3325
$ws.Range("H94").Value = 2237.1177
$ws.Range("I94").Value = 1451.2
$ws.Range("K94").Value = 1451.2
$ws.Range("M94").Value = -1000.2
$ws.Range("H99").Value = 2174.7234
$ws.Range("J99").Value = 4048.4167
$ws.Range("L99").Value = 4048.4167
$ws.Range("N99").Value = -7044.4167
$ws.Range("H105").Value = 1613.1111
$ws.Range("I105").Value = 1627.25
$ws.Range("K105").Value = 1627.25
$ws.Range("M105").Value = 119.75
$ws.Range("H126").Value = 2174.7234
$ws.Range("J126").Value = 4048.4167
$ws.Range("L126").Value = 12145.2501
$ws.Range("N126").Value = -17085.2501
$ws.Range("H132").Value = 3494.9412
$ws.Range("I132").Value = 2766.4443
$ws.Range("J132").Value = 6304.857
$ws.Range("K132").Value = 8299.332900000001
$ws.Range("L132").Value = 18914.571
$ws.Range("M132").Value = -5769.332900000001
$ws.Range("N132").Value = -23974.571
$ws.Range("H134").Value = 4137.077
$ws.Range("I134").Value = 3295.2666
$ws.Range("J134").Value = 6943.1113
$ws.Range("K134").Value = 9885.799800000001
$ws.Range("L134").Value = 20829.3339
$ws.Range("M134").Value = -7350.799800000001
$ws.Range("N134").Value = -25899.3339
$ws.Range("H136").Value = 3924.132
$ws.Range("I136").Value = 2169.8647
$ws.Range("J136").Value = 7980.875
$ws.Range("K136").Value = 6509.5941
$ws.Range("L136").Value = 23942.625
$ws.Range("M136").Value = -3959.5941
$ws.Range("N136").Value = -29042.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 118.27273
$ws.Range("I18").Value = 90.09999999999999
$ws.Range("K18").Value = 270.3
$ws.Range("M18").Value = -101.3
$ws.Range("H82").Value = 7499.5
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 7499.5
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""
$ws.Range("H115").Value = 2806.75
$ws.Range("I115").Value = 2806.75
$ws.Range("K115").Value = 8420.25
$ws.Range("M115").Value = -7245.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1224.125
$ws.Range("J2").Value = 3386.5
$ws.Range("L2").Value = 3386.5
$ws.Range("N2").Value = -3612.5
$ws.Range("H80").Value = 3246.6177
$ws.Range("I80").Value = 3111.7273
$ws.Range("J80").Value = 3311.1304
$ws.Range("K80").Value = 3111.7273
$ws.Range("L80").Value = 3311.1304
$ws.Range("M80").Value = -2113.7273
$ws.Range("N80").Value = -5307.1304
$ws.Range("H83").Value = 3246.6177
$ws.Range("I83").Value = 3111.7273
$ws.Range("J83").Value = 3311.1304
$ws.Range("K83").Value = 15558.6365
$ws.Range("L83").Value = 16555.652
$ws.Range("M83").Value = -10566.6365
$ws.Range("N83").Value = -26539.652
$ws.Range("H113").Value = 4211.4736
$ws.Range("I113").Value = 2320.8572
$ws.Range("J113").Value = 9505.200000000001
$ws.Range("K113").Value = 2320.8572
$ws.Range("L113").Value = 9505.200000000001
$ws.Range("M113").Value = -150.8571999999999
$ws.Range("N113").Value = -13845.2
$ws.Range("H122").Value = 3195.75
$ws.Range("I122").Value = 927.75
$ws.Range("K122").Value = 2783.25
$ws.Range("M122").Value = -333.25
$ws.Range("H132").Value = 4266.878
$ws.Range("I132").Value = 3705.138
$ws.Range("J132").Value = 5624.4165
$ws.Range("K132").Value = 11115.414
$ws.Range("L132").Value = 16873.2495
$ws.Range("M132").Value = -8585.414000000001
$ws.Range("N132").Value = -21933.2495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6709.613
$ws.Range("I136").Value = 5856.2085
$ws.Range("K136").Value = 17568.6255
$ws.Range("M136").Value = -15018.6255

